$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create new row 28 ------------------------------------------------
# Clone row 27's current formatting + values into row 28 first (row 27
# still carries its ORIGINAL per-column styles at this point, which is
# exactly the style pattern row 28 needs: A=14, B:H=15, I:J=15, K:N=12,
# O=13).
$ws.Range("A27:O27").Copy($ws.Range("A28:O28"))

# Blank out the columns that must stay empty on row 28.
$ws.Range("E28:H28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("N28:O28").ClearContents()

# I28/J28 need the "339018 / AUXILIO FINANCEIRO A ESTUDANTES" text
# (same text already used on row 26) while keeping row 28's own style.
# Copy-pasting VALUES ONLY from I26:J26 swaps the text in without
# disturbing the style that was just cloned from row 27.
$ws.Range("I26:J26").Copy()
$ws.Range("I28:J28").PasteSpecial(-4163)

# Numeric cells for the new row.
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 13).Value = 5641

# --- Update existing values on rows 25-27 ------------------------------
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).Value = 154996

# Row 27's Natureza Despesa cells (I27/J27) switch from the "header"
# style (s=15) to the regular style (s=11); easiest way to change just
# the style while preserving the existing shared-string text is to
# paste the formats from a cell that already carries style 11 (I25/J25).
$ws.Range("I25:J25").Copy()
$ws.Range("I27:J27").PasteSpecial(-4122)
$ws.Cells.Item(27, 11).Value = 3354

# --- Extend the merged blocks down to row 28 ---------------------------
$ws.Range("A25:A28").Merge()
$ws.Range("B25:B28").Merge()
$ws.Range("C25:C28").Merge()
$ws.Range("D25:D28").Merge()
$ws.Range("E27:E28").Merge()
$ws.Range("F27:F28").Merge()
$ws.Range("G27:G28").Merge()
$ws.Range("H27:H28").Merge()
